$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix spelling mistake: "Test Box Interactions" -> "Text Box Interactions"
$a2 = $ws.Range("A2")
$a2.Value = "Text Box Interactions"

# Normalize C2's formatting so it no longer needs its own distinct
# "applyFill" cell-format record (matches the rest of row 2, e.g. D2).
$c2 = $ws.Range("C2")
$c2.Interior.Pattern = [Microsoft.Office.Interop.Excel.XlPattern]::xlPatternNone

# Update the saved selection/active cell.
$ws.Range("A4").Select() | Out-Null
